$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7708112001419067
$ws.Range("B1").Value = 1.37587559223175
$ws.Range("C1").Value = 4.113471031188965
$ws.Range("D1").Value = 5.913780689239502
$ws.Range("E1").Value = 1.742050528526306
